# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets, which contain duplicate data tables.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    3  = 3099
    5  = 2666
    7  = 140
    9  = 1407
    11 = 62
    13 = 1201
    16 = 328
    17 = 38
    18 = 35
    22 = 2577
    23 = 39
    24 = 294
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
